$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D, J, K, L, M, P across rows 2-6
# (derived from a cyclic rotation of the original row data)
$data = @{
    2 = @{ D = 44175; J = 1400; K = 1900; L = 2000; M = 1950; P = 1950 }
    3 = @{ D = 44638; J = 800;  K = 2500; L = 2800; M = 2650; P = 2650 }
    4 = @{ D = 44200; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
    5 = @{ D = 44210; J = 1450; K = 1600; L = 1700; M = 1650; P = 1650 }
    6 = @{ D = 44537; J = 800;  K = 1300; L = 1400; M = 1350; P = 1350 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value2 = $vals.D
    $ws.Range("J$row").Value2 = $vals.J
    $ws.Range("K$row").Value2 = $vals.K
    $ws.Range("L$row").Value2 = $vals.L
    $ws.Range("M$row").Value2 = $vals.M
    $ws.Range("P$row").Value2 = $vals.P
}
